# Underline the whole first paragraph ("Prova ProVa2" + its paragraph
# mark) and re-anchor the existing "_GoBack" bookmark so that it spans
# from the start of that paragraph through to the start of the next
# (empty) paragraph.

$d = $word.ActiveDocument

# --- 1. Re-position the "_GoBack" bookmark --------------------------------
# It currently sits collapsed between the runs "V" and "a2". Move it so it
# wraps the entire first paragraph, starting right before "Prova" and
# ending at the very start of the second (empty) paragraph.
$firstPara = $d.Paragraphs.Item(1)
$bookmarkRange = $d.Range($firstPara.Range.Start, $d.Content.End)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 2. Underline the whole paragraph (text + paragraph mark) -------------
$firstPara.Range.Font.Underline = 1
